$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New confusion-table data (intent, confused_with, count) for rows 2..49,
# reflecting additional NLU training examples that shifted the confusion
# pairs / counts captured by the failed_stories report.
$data = @(
  @("comment_positive", "vocative_thank_you", 13),
  @("vocative_thank_you", "comment_positive", 8),
  @("features_date", "user_love", 4),
  @("covid_worry", "user_scared", 3),
  @("covid_origins", "covid_info", 3),
  @("bot_capabilities", "vocative_call", 3),
  @("covid_dangerous", "spread_risk", 3),
  @("bot_personal_questions", "bot_residence", 3),
  @("bot_capabilities", "vocative_help", 2),
  @("prevention_home", "prevention_medical_attention", 2),
  @("bot_residence", "bot_origin", 2),
  @("comment_positive", "comment_smart", 2),
  @("bot_sexual", "bot_personal_questions", 2),
  @("vocative_call", "cc_lets_talk", 2),
  @("covid_info", "covid_origins", 2),
  @("vocative_call", "user_friend", 2),
  @("covid_duration", "covid_crisis_howlong", 1),
  @("covid_duration", "covid_info", 1),
  @("vocative_yes", "country", 1),
  @("vocative_no", "comment_negative", 1),
  @("greeting_hello", "start", 1),
  @("quaratine_how_it_works", "quarantine_general", 1),
  @("prevention_home", "comment_negative", 1),
  @("covid_situation_last_update", "country", 1),
  @("bot_name", "bot_capabilities", 1),
  @("bot_name", "bot_personal_questions", 1),
  @("covid_crisis_howlong", "spread_general", 1),
  @("covid_symptoms", "test_virus", 1),
  @("bot_personality", "bot_capabilities", 1),
  @("spread_risk", "spread_general", 1),
  @("comment_racist", "comment_offense", 1),
  @("comment_racist", "cc_geography", 1),
  @("bot_sing", "bot_personal_questions", 1),
  @("covid_treatment", "prevention_medicine", 1),
  @("covid_info", "coronavirus_info", 1),
  @("covid_origins", "cc_philosophical", 1),
  @("user_particles", "comment_smart", 1),
  @("user_particles", "vocative_you_welcome", 1),
  @("comment_negative", "comment_offense", 1),
  @("comment_negative", "bot_capabilities", 1),
  @("bot_personal_questions", "user_love", 1),
  @("prevention_medical_attention", "covid_procedure_after_infection", 1),
  @("prevention_medical_attention", "covid_symptoms", 1),
  @("cc_philosophical", "features_time", 1),
  @("cc_philosophical", "user_love", 1),
  @("bot_music", "bot_sing", 1),
  @("test_who", "test_virus", 1),
  @("covid_symptoms", "prevention_medical_attention", 1)
)

$row = 2
foreach ($entry in $data) {
  $ws.Cells.Item($row, 1).Value = $entry[0]
  $ws.Cells.Item($row, 2).Value = $entry[1]
  $ws.Cells.Item($row, 3).Value = $entry[2]
  $row = $row + 1
}

# The old table had 10 more rows (50-59); remove that now-unused tail so the
# sheet's used range / dimension shrinks back to A1:C49.
$ws.Range("A50:C59").ClearContents()
